# Sample Project / Main.xlsx — "Rules" sheet
# The R40 rule's "From" value (cell B11) is retyped as the text "1".
# Leading the literal with an apostrophe forces Excel to store it as text
# (a new shared-string entry) instead of re-interpreting it as a number,
# matching the original cell's text ("s") type.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("B11").Value = "'1"
